$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns before the existing ExpPoints column (C),
# shifting ExpPoints from C to H.
$ws.Range("C1:G1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"

# Updated ExpPoints values (now in column H) - matchday 11 prediction.
$ws.Range("H2").Value = 89.40649984141315
$ws.Range("H3").Value = 88.96839385999897
$ws.Range("H4").Value = 71.00687961084466
$ws.Range("H5").Value = 65.33208261035935
$ws.Range("H6").Value = 61.01722997278711
$ws.Range("H7").Value = 58.54224499439288
$ws.Range("H8").Value = 55.32139877381939
$ws.Range("H9").Value = 49.24423320919442
$ws.Range("H10").Value = 49.03313109624208
$ws.Range("H11").Value = 48.87741563842056
$ws.Range("H12").Value = 48.19963321355022
$ws.Range("H13").Value = 47.41724244679676
$ws.Range("H14").Value = 46.68019926713524
$ws.Range("H15").Value = 46.55388650526439
$ws.Range("H16").Value = 39.27877032899355
$ws.Range("H17").Value = 37.15377576108096
$ws.Range("H18").Value = 35.81359384052848
$ws.Range("H19").Value = 34.36662415321098
$ws.Range("H20").Value = 33.03814174691885
$ws.Range("H21").Value = 28.39376837986908

# Re-order a few teams (table shuffled slightly for matchday 11).
$ws.Range("B9").Value = "Real Sociedad"
$ws.Range("B10").Value = "Espanyol"
$ws.Range("B11").Value = "Osasuna"

$ws.Range("B13").Value = "Sevilla"
$ws.Range("B14").Value = "Getafe"
$ws.Range("B15").Value = "Celta de Vigo"
